$d = $word.ActiveDocument

function Fix-FirstLetter($paraIndex, $newFirstChar, $splitIntoTwoRuns) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range

    # Replace just the first character of the caption text.
    $first = $d.Range($r.Start, $r.Start + 1)
    $first.Text = $newFirstChar

    if ($splitIntoTwoRuns) {
        # Force the first character to live in its own run (distinct from
        # the rest of the caption) by toggling a character formatting
        # property on it and then reverting it. Re-fetch the range each
        # time since the previous call can reseat it.
        $a1 = $d.Range($r.Start, $r.Start + 1)
        $a1.Font.Bold = $true
        $a2 = $d.Range($r.Start, $r.Start + 1)
        $a2.Font.Bold = $false
    }

    # Re-fetch the paragraph range (text length may have changed) and make
    # sure the remaining caption run(s) keep an explicit (even if empty)
    # run-properties element, matching the saved document's run shape.
    $p2 = $d.Paragraphs.Item($paraIndex)
    $r2 = $p2.Range
    if ($splitIntoTwoRuns) {
        $restStart = $r2.Start + 1
    } else {
        $restStart = $r2.Start
    }
    $rest = $d.Range($restStart, $r2.End - 1)
    if ($rest.Start -lt $rest.End) {
        $rest.Font.Bold = $true
        $rest2 = $d.Range($restStart, $r2.End - 1)
        $rest2.Font.Bold = $false
    }
}

Fix-FirstLetter 6 "р" $true
Fix-FirstLetter 9 "р" $true
Fix-FirstLetter 12 "р" $false
